$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Row 60: fill in the "status" column (G60) which was left blank before.
# ---------------------------------------------------------------------------
$ws.Range("G60").Value = "Done!"

# ---------------------------------------------------------------------------
# Helper-ish: style a new data row the same way the existing rows 56-61 are
# styled (green fill, thin border all round, column A/H/I/K word-wrapped,
# column L formatted as a date). Column G and J are intentionally left
# untouched (no fill / no border) to match the existing rows that have no
# status / no J value.
# ---------------------------------------------------------------------------

function Format-ExperimentRow($rowNum) {
    $greenRanges = @("A$rowNum`:F$rowNum", "H$rowNum`:I$rowNum", "K$rowNum`:L$rowNum")
    foreach ($addr in $greenRanges) {
        $rng = $ws.Range($addr)
        $rng.Interior.Color = 5287936
        $rng.Borders.LineStyle = 1
    }
    $ws.Range("A$rowNum").WrapText = $true
    $ws.Range("H$rowNum").WrapText = $true
    $ws.Range("I$rowNum").WrapText = $true
    $ws.Range("K$rowNum").WrapText = $true
    $ws.Range("L$rowNum").NumberFormat = "d-mmm-yy"
}

# ---------------------------------------------------------------------------
# 2) Row 62 - EXP60 (run_id 53)
# ---------------------------------------------------------------------------
Format-ExperimentRow 62

$ws.Range("A62").Value = 'TPR concatenated with LSTM in 
phrase embedding layer 
batchsize = 40. With visualizations. With regularization. Regularization weights=0.00002. nRoles=20, nSymbols=100 [running from "QA_TPR_for_Run" branch "master" ]. '
$ws.Range("A62").Characters(1, 26).Font.Bold = $true
$ws.Range("A62").Characters(27, 87).Font.Bold = $false
$ws.Range("A62").Characters(114, 55).Font.Bold = $true
$ws.Range("A62").Characters(169, 51).Font.Bold = $false

$ws.Range("D62").Value = "EXP60.txt"

$ws.Range("B62").Value = 'python -m basic.cli --mode train --noload --len_opt --cluster --LSTMandTPR True --TPRregularizer1 True --TPRvis True --cF 0.00002 --cR 0.00002 --nRoles 20 --nSymbols 100 --batch_size 40 --run_id 53 |& tee /home/hpalangi/QA/TPR_Stuff/Codes/TPR_ver1.0/Log_Files/EXP60.txt'

$ws.Range("C62").Value = "DLDGX / 0"
$ws.Range("E62").Value = 53
$ws.Range("F62").Value = 2

$ws.Rows("62").RowHeight = 180

# ---------------------------------------------------------------------------
# 3) Row 63 - EXP61 (run_id 54)
# ---------------------------------------------------------------------------
Format-ExperimentRow 63

$ws.Range("D63").Value = "EXP61.txt"

$ws.Range("A63").Value = 'TPR concatenated with LSTM in 
phrase embedding layer 
batchsize = 40. With visualizations. With regularization. Regularization weights=0.00003. nRoles=20, nSymbols=100 [running from "QA_TPR_for_Run" branch "master" ]. '
$ws.Range("A63").Characters(1, 26).Font.Bold = $true
$ws.Range("A63").Characters(27, 87).Font.Bold = $false
$ws.Range("A63").Characters(114, 55).Font.Bold = $true
$ws.Range("A63").Characters(169, 51).Font.Bold = $false

$ws.Range("B63").Value = 'python -m basic.cli --mode train --noload --len_opt --cluster --LSTMandTPR True --TPRregularizer1 True --TPRvis True --cF 0.00003 --cR 0.00003 --nRoles 20 --nSymbols 100 --batch_size 40 --run_id 54 |& tee /home/hpalangi/QA/TPR_Stuff/Codes/TPR_ver1.0/Log_Files/EXP61.txt'

$ws.Range("C63").Value = "DLDGX / 4"
$ws.Range("E63").Value = 54
$ws.Range("F63").Value = 4

$ws.Rows("63").RowHeight = 180

# ---------------------------------------------------------------------------
# 4) Row 64 - EXP62 (run_id 55)
# ---------------------------------------------------------------------------
Format-ExperimentRow 64

$ws.Range("A64").Value = 'TPR concatenated with LSTM in 
phrase embedding layer 
batchsize = 60. With visualizations. With regularization. Regularization weights=0.00003. dRoles=5, dSymbols=5, nRoles=20, nSymbols=100 [running from "QA_TPR_for_Run" branch "master" ]. '
$ws.Range("A64").Characters(1, 26).Font.Bold = $true
$ws.Range("A64").Characters(27, 29).Font.Bold = $false
$ws.Range("A64").Characters(56, 14).Font.Bold = $true
$ws.Range("A64").Characters(70, 44).Font.Bold = $false
$ws.Range("A64").Characters(114, 77).Font.Bold = $true
$ws.Range("A64").Characters(191, 51).Font.Bold = $false

$ws.Range("D64").Value = "EXP62.txt"

$ws.Range("B64").Value = 'python -m basic.cli --mode train --noload --len_opt --cluster --LSTMandTPR True --TPRregularizer1 True --TPRvis True --cF 0.00003 --cR 0.00003 --nRoles 20 --nSymbols 100 --dRoles 5 --dSymbols 5 --batch_size 60 --run_id 55 |& tee /home/hpalangi/QA/TPR_Stuff/Codes/TPR_ver1.0/Log_Files/EXP62.txt'

$ws.Range("C64").Value = "DLDGX / 5"
$ws.Range("E64").Value = 55
$ws.Range("F64").Value = 5

$ws.Rows("64").RowHeight = 195

# ---------------------------------------------------------------------------
# 5) Move the selection / freeze-pane viewport down to the newly added rows.
# ---------------------------------------------------------------------------
$ws.Range("B64").Select()
